$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (existing C/D/E/F shift right to D/E/F/G)
$ws.Columns("C").Insert()

# Fill the new column C with the ORCID identifiers, row by row (rows 1-22)
$orcids = @(
    "0000-0003-3388-4343",
    "0000-0003-3388-4344",
    "0000-0003-3388-4345",
    "0000-0003-3388-4346",
    "0000-0003-3388-4347",
    "0000-0003-3388-4348",
    "0000-0003-3388-4349",
    "0000-0003-3388-4350",
    "0000-0003-3388-4351",
    "0000-0003-3388-4352",
    "0000-0003-3388-4353",
    "0000-0003-3388-4354",
    "0000-0003-3388-4355",
    "0000-0003-3388-4356",
    "0000-0003-3388-4357",
    "0000-0003-3388-4358",
    "0000-0003-3388-4359",
    "0000-0003-3388-4360",
    "0000-0003-3388-4361",
    "0000-0003-3388-4362",
    "0000-0003-3388-4363",
    "0000-0003-3388-4364"
)

for ($i = 0; $i -lt $orcids.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $orcids[$i]
}

# Match column C's width to column B's width (same visual width as names column)
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Update the active selection to match the post-edit state (single cell C1)
$ws.Range("C1").Select()
